$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.646.72"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.120.67"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.27%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +1.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.92"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5256"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.71%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4554"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.58"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09109"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.176"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.49"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.59%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.127.32"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.863"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.149"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +5.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001174"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "97.21"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.29%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06693"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.47"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.011"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.315"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.715.31"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.93"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +5.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.365"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.365.25"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.38"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.02"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.555"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.94"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.207"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.83%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.647"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.364"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.27%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.69"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.852"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +7.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02637"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06876"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2330"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.63"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6903"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.260"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.91"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.61%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6488"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.315"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000372"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +22.58%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.256"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.88%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.45"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.196"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.30%  "
